$d = $word.ActiveDocument

# The existing "_GoBack" bookmark sits at the end of the "Background
# Image" bullet. It needs to move to the new bullet we are about to
# add (right after "Fire sprite on fuel bar"), so drop it here; it is
# recreated below at the correct spot.
$d.Bookmarks("_GoBack").Delete()

# Add a new list-item paragraph right after the "Background Image"
# bullet; InsertParagraphAfter() inherits the ListParagraph style /
# numbering (numId 1, ilvl 0) and the run language formatting.
$srcPara = $d.Paragraphs(3)
$srcPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(4)

# Type the whole line as a single run first ...
$line = "Fire sprite on fuel bar: Everett Ellerbrock"
$newPara.Range.InsertAfter($line)

$paraStart = $newPara.Range.Start
$afterFirstSentence = $paraStart + ("Fire sprite on fuel bar".Length)
$afterColonSpace = $afterFirstSentence + 2

# ... then split it in two places using collapsed bookmarks: first
# recreate "_GoBack" between "Fire sprite on fuel bar" and ": ", then
# a transient bookmark between ": " and "Everett Ellerbrock" so that
# text keeps living in three separate runs once saved.
$gbRange = $d.Range($afterFirstSentence, $afterFirstSentence)
$d.Bookmarks.Add("_GoBack", $gbRange)

$splitRange = $d.Range($afterColonSpace, $afterColonSpace)
$d.Bookmarks.Add("TempRunSplit", $splitRange)
$d.Bookmarks("TempRunSplit").Delete()
